{"js": "// Replace each arithmetic-expression cell in the (only) table with its\n// updated expression. Cells are addressed by zero-based (row, col) so the\n// edit is applied positionally -- this matters because a couple of old\n// expressions (e.g. \"98-55=\") repeat in the sheet but map to different\n// replacements, so a blind global text replace would be wrong.\nconst replacements = [\n  { r: 0, c: 0, oldText: \"86-81=\", newText: \"51+2=\" },\n  { r: 0, c: 1, oldText: \"12+33=\", newText: \"57+29=\" },\n  { r: 0, c: 2, oldText: \"78-24=\", newText: \"20+16=\" },\n  { r: 0, c: 3, oldText: \"72-14=\", newText: \"14+24=\" },\n  { r: 0, c: 4, oldText: \"29+7=\", newText: \"93-10=\" },\n  { r: 1, c: 0, oldText: \"79+2=\", newText: \"20+46=\" },\n  { r: 1, c: 1, oldText: \"28-26=\", newText: \"5+29=\" },\n  { r: 1, c: 2, oldText: \"68+13=\", newText: \"67-50=\" },\n  { r: 1, c: 3, oldText: \"48+7=\", newText: \"48-43=\" },\n  { r: 1, c: 4, oldText: \"6+78=\", newText: \"89-34=\" },\n  { r: 2, c: 0, oldText: \"55-24=\", newText: \"45-19=\" },\n  { r: 2, c: 1, oldText: \"24+57=\", newText: \"51-22=\" },\n  { r: 2, c: 2, oldText: \"13+83=\", newText: \"16+47=\" },\n  { r: 2, c: 3, oldText: \"88-80=\", newText: \"8+89=\" },\n  { r: 2, c: 4, oldText: \"98-50=\", newText: \"45+5=\" },\n  { r: 3, c: 0, oldText: \"55+13=\", newText: \"31-30=\" },\n  { r: 3, c: 1, oldText: \"19+1=\", newText: \"79-49=\" },\n  { r: 3, c: 2, oldText: \"9+6=\", newText: \"0+92=\" },\n  { r: 3, c: 3, oldText: \"20+41=\", newText: \"89-73=\" },\n  { r: 3, c: 4, oldText: \"58+21=\", newText: \"30-4=\" },\n  { r: 4, c: 0, oldText: \"29-9=\", newText: \"63-33=\" },\n  { r: 4, c: 1, oldText: \"93-40=\", newText: \"41-9=\" },\n  { r: 4, c: 2, oldText: \"55-25=\", newText: \"46+37=\" },\n  { r: 4, c: 3, oldText: \"15+32=\", newText: \"65+0=\" },\n  { r: 4, c: 4, oldText: \"36+33=\", newText: \"44-9=\" },\n  { r: 5, c: 0, oldText: \"36+63=\", newText: \"55+30=\" },\n  { r: 5, c: 1, oldText: \"49-45=\", newText: \"46+21=\" },\n  { r: 5, c: 2, oldText: \"18+9=\", newText: \"44+53=\" },\n  { r: 5, c: 3, oldText: \"74-48=\", newText: \"34-29=\" },\n  { r: 5, c: 4, oldText: \"59+1=\", newText: \"11-10=\" },\n  { r: 6, c: 0, oldText: \"79-69=\", newText: \"68-51=\" },\n  { r: 6, c: 1, oldText: \"13+44=\", newText: \"33+4=\" },\n  { r: 6, c: 2, oldText: \"28+40=\", newText: \"4+10=\" },\n  { r: 6, c: 3, oldText: \"7+88=\", newText: \"29+27=\" },\n  { r: 6, c: 4, oldText: \"8+11=\", newText: \"1+65=\" },\n  { r: 7, c: 0, oldText: \"57-4=\", newText: \"9+30=\" },\n  { r: 7, c: 1, oldText: \"90-35=\", newText: \"59-12=\" },\n  { r: 7, c: 2, oldText: \"13+78=\", newText: \"3+50=\" },\n  { r: 7, c: 3, oldText: \"0+50=\", newText: \"76-18=\" },\n  { r: 7, c: 4, oldText: \"6+2=\", newText: \"99-20=\" },\n  { r: 8, c: 0, oldText: \"43+56=\", newText: \"28-10=\" },\n  { r: 8, c: 1, oldText: \"12+84=\", newText: \"2+23=\" },\n  { r: 8, c: 2, oldText: \"41-17=\", newText: \"19+8=\" },\n  { r: 8, c: 3, oldText: \"32+5=\", newText: \"19+24=\" },\n  { r: 8, c: 4, oldText: \"98-55=\", newText: \"60+23=\" },\n  { r: 9, c: 0, oldText: \"97-59=\", newText: \"52+9=\" },\n  { r: 9, c: 1, oldText: \"0+77=\", newText: \"49+27=\" },\n  { r: 9, c: 2, oldText: \"65-22=\", newText: \"51+27=\" },\n  { r: 9, c: 3, oldText: \"28+46=\", newText: \"12+48=\" },\n  { r: 9, c: 4, oldText: \"3+89=\", newText: \"68-7=\" },\n  { r: 10, c: 0, oldText: \"65-59=\", newText: \"41-9=\" },\n  { r: 10, c: 1, oldText: \"81-31=\", newText: \"20+17=\" },\n  { r: 10, c: 2, oldText: \"65-49=\", newText: \"52+6=\" },\n  { r: 10, c: 3, oldText: \"87-36=\", newText: \"23+21=\" },\n  { r: 10, c: 4, oldText: \"70+28=\", newText: \"67-14=\" },\n  { r: 11, c: 0, oldText: \"98-59=\", newText: \"79-13=\" },\n  { r: 11, c: 1, oldText: \"25-10=\", newText: \"82-37=\" },\n  { r: 11, c: 2, oldText: \"87-2=\", newText: \"76-30=\" },\n  { r: 11, c: 3, oldText: \"94-76=\", newText: \"70-32=\" },\n  { r: 11, c: 4, oldText: \"17-12=\", newText: \"50-42=\" },\n  { r: 12, c: 0, oldText: \"87-11=\", newText: \"76-46=\" },\n  { r: 12, c: 1, oldText: \"66+19=\", newText: \"60+3=\" },\n  { r: 12, c: 2, oldText: \"18+53=\", newText: \"0+5=\" },\n  { r: 12, c: 3, oldText: \"30+10=\", newText: \"71-52=\" },\n  { r: 12, c: 4, oldText: \"15-10=\", newText: \"9+90=\" },\n  { r: 13, c: 0, oldText: \"76+1=\", newText: \"71+27=\" },\n  { r: 13, c: 1, oldText: \"69-6=\", newText: \"34+24=\" },\n  { r: 13, c: 2, oldText: \"82-78=\", newText: \"98-58=\" },\n  { r: 13, c: 3, oldText: \"31+50=\", newText: \"41+12=\" },\n  { r: 13, c: 4, oldText: \"66-45=\", newText: \"50+45=\" },\n  { r: 14, c: 0, oldText: \"74+13=\", newText: \"19+34=\" },\n  { r: 14, c: 1, oldText: \"75-53=\", newText: \"44-29=\" },\n  { r: 14, c: 2, oldText: \"72-71=\", newText: \"82-30=\" },\n  { r: 14, c: 3, oldText: \"57+31=\", newText: \"32+31=\" },\n  { r: 14, c: 4, oldText: \"27+12=\", newText: \"22+56=\" },\n  { r: 15, c: 0, oldText: \"37+5=\", newText: \"45-11=\" },\n  { r: 15, c: 1, oldText: \"43-25=\", newText: \"92-4=\" },\n  { r: 15, c: 2, oldText: \"73-50=\", newText: \"50-7=\" },\n  { r: 15, c: 3, oldText: \"10-8=\", newText: \"70-9=\" },\n  { r: 15, c: 4, oldText: \"58-43=\", newText: \"88-20=\" },\n  { r: 16, c: 0, oldText: \"98-55=\", newText: \"90-50=\" },\n  { r: 16, c: 1, oldText: \"19+19=\", newText: \"71-56=\" },\n  { r: 16, c: 2, oldText: \"1+21=\", newText: \"27+48=\" },\n  { r: 16, c: 3, oldText: \"99-18=\", newText: \"13+69=\" },\n  { r: 16, c: 4, oldText: \"1+50=\", newText: \"72-24=\" },\n  { r: 17, c: 0, oldText: \"95-71=\", newText: \"17+6=\" },\n  { r: 17, c: 1, oldText: \"21+9=\", newText: \"44+54=\" },\n  { r: 17, c: 2, oldText: \"49+24=\", newText: \"81-55=\" },\n  { r: 17, c: 3, oldText: \"79-22=\", newText: \"71-55=\" },\n  { r: 17, c: 4, oldText: \"80-79=\", newText: \"53+35=\" },\n  { r: 18, c: 0, oldText: \"55+31=\", newText: \"88-62=\" },\n  { r: 18, c: 1, oldText: \"62+16=\", newText: \"50+39=\" },\n  { r: 18, c: 2, oldText: \"63-43=\", newText: \"90-48=\" },\n  { r: 18, c: 3, oldText: \"92-77=\", newText: \"83-22=\" },\n  { r: 18, c: 4, oldText: \"83-65=\", newText: \"62+23=\" },\n  { r: 19, c: 0, oldText: \"86-45=\", newText: \"32+2=\" },\n  { r: 19, c: 1, oldText: \"11+6=\", newText: \"69+17=\" },\n  { r: 19, c: 2, oldText: \"14+48=\", newText: \"49+44=\" },\n  { r: 19, c: 3, oldText: \"14+14=\", newText: \"76-56=\" },\n  { r: 19, c: 4, oldText: \"14+54=\", newText: \"58+25=\" }\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\n\n// Load current cell values first so we can sanity-check them against the\n// expected \"before\" text (purely informational; the edit itself is applied\n// positionally regardless, matching how the source diff addresses cells).\nconst cells = replacements.map((p) => table.getCell(p.r, p.c));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const entry = replacements[i];\n  cells[i].value = entry.newText;\n}\n\nawait context.sync();\n", "ps1": "# Replace each arithmetic-expression cell in the (only) table with its\n# updated expression. Cells are addressed by 1-based (Row, Col) so the\n# edit is applied positionally -- this matters because a couple of old\n# expressions (e.g. \"98-55=\") repeat in the sheet but map to different\n# replacements, so a blind global Find/Replace by text would be wrong.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"86-81=\"; NewText = \"51+2=\" },\n    @{ Row = 1; Col = 2; OldText = \"12+33=\"; NewText = \"57+29=\" },\n    @{ Row = 1; Col = 3; OldText = \"78-24=\"; NewText = \"20+16=\" },\n    @{ Row = 1; Col = 4; OldText = \"72-14=\"; NewText = \"14+24=\" },\n    @{ Row = 1; Col = 5; OldText = \"29+7=\"; NewText = \"93-10=\" },\n    @{ Row = 2; Col = 1; OldText = \"79+2=\"; NewText = \"20+46=\" },\n    @{ Row = 2; Col = 2; OldText = \"28-26=\"; NewText = \"5+29=\" },\n    @{ Row = 2; Col = 3; OldText = \"68+13=\"; NewText = \"67-50=\" },\n    @{ Row = 2; Col = 4; OldText = \"48+7=\"; NewText = \"48-43=\" },\n    @{ Row = 2; Col = 5; OldText = \"6+78=\"; NewText = \"89-34=\" },\n    @{ Row = 3; Col = 1; OldText = \"55-24=\"; NewText = \"45-19=\" },\n    @{ Row = 3; Col = 2; OldText = \"24+57=\"; NewText = \"51-22=\" },\n    @{ Row = 3; Col = 3; OldText = \"13+83=\"; NewText = \"16+47=\" },\n    @{ Row = 3; Col = 4; OldText = \"88-80=\"; NewText = \"8+89=\" },\n    @{ Row = 3; Col = 5; OldText = \"98-50=\"; NewText = \"45+5=\" },\n    @{ Row = 4; Col = 1; OldText = \"55+13=\"; NewText = \"31-30=\" },\n    @{ Row = 4; Col = 2; OldText = \"19+1=\"; NewText = \"79-49=\" },\n    @{ Row = 4; Col = 3; OldText = \"9+6=\"; NewText = \"0+92=\" },\n    @{ Row = 4; Col = 4; OldText = \"20+41=\"; NewText = \"89-73=\" },\n    @{ Row = 4; Col = 5; OldText = \"58+21=\"; NewText = \"30-4=\" },\n    @{ Row = 5; Col = 1; OldText = \"29-9=\"; NewText = \"63-33=\" },\n    @{ Row = 5; Col = 2; OldText = \"93-40=\"; NewText = \"41-9=\" },\n    @{ Row = 5; Col = 3; OldText = \"55-25=\"; NewText = \"46+37=\" },\n    @{ Row = 5; Col = 4; OldText = \"15+32=\"; NewText = \"65+0=\" },\n    @{ Row = 5; Col = 5; OldText = \"36+33=\"; NewText = \"44-9=\" },\n    @{ Row = 6; Col = 1; OldText = \"36+63=\"; NewText = \"55+30=\" },\n    @{ Row = 6; Col = 2; OldText = \"49-45=\"; NewText = \"46+21=\" },\n    @{ Row = 6; Col = 3; OldText = \"18+9=\"; NewText = \"44+53=\" },\n    @{ Row = 6; Col = 4; OldText = \"74-48=\"; NewText = \"34-29=\" },\n    @{ Row = 6; Col = 5; OldText = \"59+1=\"; NewText = \"11-10=\" },\n    @{ Row = 7; Col = 1; OldText = \"79-69=\"; NewText = \"68-51=\" },\n    @{ Row = 7; Col = 2; OldText = \"13+44=\"; NewText = \"33+4=\" },\n    @{ Row = 7; Col = 3; OldText = \"28+40=\"; NewText = \"4+10=\" },\n    @{ Row = 7; Col = 4; OldText = \"7+88=\"; NewText = \"29+27=\" },\n    @{ Row = 7; Col = 5; OldText = \"8+11=\"; NewText = \"1+65=\" },\n    @{ Row = 8; Col = 1; OldText = \"57-4=\"; NewText = \"9+30=\" },\n    @{ Row = 8; Col = 2; OldText = \"90-35=\"; NewText = \"59-12=\" },\n    @{ Row = 8; Col = 3; OldText = \"13+78=\"; NewText = \"3+50=\" },\n    @{ Row = 8; Col = 4; OldText = \"0+50=\"; NewText = \"76-18=\" },\n    @{ Row = 8; Col = 5; OldText = \"6+2=\"; NewText = \"99-20=\" },\n    @{ Row = 9; Col = 1; OldText = \"43+56=\"; NewText = \"28-10=\" },\n    @{ Row = 9; Col = 2; OldText = \"12+84=\"; NewText = \"2+23=\" },\n    @{ Row = 9; Col = 3; OldText = \"41-17=\"; NewText = \"19+8=\" },\n    @{ Row = 9; Col = 4; OldText = \"32+5=\"; NewText = \"19+24=\" },\n    @{ Row = 9; Col = 5; OldText = \"98-55=\"; NewText = \"60+23=\" },\n    @{ Row = 10; Col = 1; OldText = \"97-59=\"; NewText = \"52+9=\" },\n    @{ Row = 10; Col = 2; OldText = \"0+77=\"; NewText = \"49+27=\" },\n    @{ Row = 10; Col = 3; OldText = \"65-22=\"; NewText = \"51+27=\" },\n    @{ Row = 10; Col = 4; OldText = \"28+46=\"; NewText = \"12+48=\" },\n    @{ Row = 10; Col = 5; OldText = \"3+89=\"; NewText = \"68-7=\" },\n    @{ Row = 11; Col = 1; OldText = \"65-59=\"; NewText = \"41-9=\" },\n    @{ Row = 11; Col = 2; OldText = \"81-31=\"; NewText = \"20+17=\" },\n    @{ Row = 11; Col = 3; OldText = \"65-49=\"; NewText = \"52+6=\" },\n    @{ Row = 11; Col = 4; OldText = \"87-36=\"; NewText = \"23+21=\" },\n    @{ Row = 11; Col = 5; OldText = \"70+28=\"; NewText = \"67-14=\" },\n    @{ Row = 12; Col = 1; OldText = \"98-59=\"; NewText = \"79-13=\" },\n    @{ Row = 12; Col = 2; OldText = \"25-10=\"; NewText = \"82-37=\" },\n    @{ Row = 12; Col = 3; OldText = \"87-2=\"; NewText = \"76-30=\" },\n    @{ Row = 12; Col = 4; OldText = \"94-76=\"; NewText = \"70-32=\" },\n    @{ Row = 12; Col = 5; OldText = \"17-12=\"; NewText = \"50-42=\" },\n    @{ Row = 13; Col = 1; OldText = \"87-11=\"; NewText = \"76-46=\" },\n    @{ Row = 13; Col = 2; OldText = \"66+19=\"; NewText = \"60+3=\" },\n    @{ Row = 13; Col = 3; OldText = \"18+53=\"; NewText = \"0+5=\" },\n    @{ Row = 13; Col = 4; OldText = \"30+10=\"; NewText = \"71-52=\" },\n    @{ Row = 13; Col = 5; OldText = \"15-10=\"; NewText = \"9+90=\" },\n    @{ Row = 14; Col = 1; OldText = \"76+1=\"; NewText = \"71+27=\" },\n    @{ Row = 14; Col = 2; OldText = \"69-6=\"; NewText = \"34+24=\" },\n    @{ Row = 14; Col = 3; OldText = \"82-78=\"; NewText = \"98-58=\" },\n    @{ Row = 14; Col = 4; OldText = \"31+50=\"; NewText = \"41+12=\" },\n    @{ Row = 14; Col = 5; OldText = \"66-45=\"; NewText = \"50+45=\" },\n    @{ Row = 15; Col = 1; OldText = \"74+13=\"; NewText = \"19+34=\" },\n    @{ Row = 15; Col = 2; OldText = \"75-53=\"; NewText = \"44-29=\" },\n    @{ Row = 15; Col = 3; OldText = \"72-71=\"; NewText = \"82-30=\" },\n    @{ Row = 15; Col = 4; OldText = \"57+31=\"; NewText = \"32+31=\" },\n    @{ Row = 15; Col = 5; OldText = \"27+12=\"; NewText = \"22+56=\" },\n    @{ Row = 16; Col = 1; OldText = \"37+5=\"; NewText = \"45-11=\" },\n    @{ Row = 16; Col = 2; OldText = \"43-25=\"; NewText = \"92-4=\" },\n    @{ Row = 16; Col = 3; OldText = \"73-50=\"; NewText = \"50-7=\" },\n    @{ Row = 16; Col = 4; OldText = \"10-8=\"; NewText = \"70-9=\" },\n    @{ Row = 16; Col = 5; OldText = \"58-43=\"; NewText = \"88-20=\" },\n    @{ Row = 17; Col = 1; OldText = \"98-55=\"; NewText = \"90-50=\" },\n    @{ Row = 17; Col = 2; OldText = \"19+19=\"; NewText = \"71-56=\" },\n    @{ Row = 17; Col = 3; OldText = \"1+21=\"; NewText = \"27+48=\" },\n    @{ Row = 17; Col = 4; OldText = \"99-18=\"; NewText = \"13+69=\" },\n    @{ Row = 17; Col = 5; OldText = \"1+50=\"; NewText = \"72-24=\" },\n    @{ Row = 18; Col = 1; OldText = \"95-71=\"; NewText = \"17+6=\" },\n    @{ Row = 18; Col = 2; OldText = \"21+9=\"; NewText = \"44+54=\" },\n    @{ Row = 18; Col = 3; OldText = \"49+24=\"; NewText = \"81-55=\" },\n    @{ Row = 18; Col = 4; OldText = \"79-22=\"; NewText = \"71-55=\" },\n    @{ Row = 18; Col = 5; OldText = \"80-79=\"; NewText = \"53+35=\" },\n    @{ Row = 19; Col = 1; OldText = \"55+31=\"; NewText = \"88-62=\" },\n    @{ Row = 19; Col = 2; OldText = \"62+16=\"; NewText = \"50+39=\" },\n    @{ Row = 19; Col = 3; OldText = \"63-43=\"; NewText = \"90-48=\" },\n    @{ Row = 19; Col = 4; OldText = \"92-77=\"; NewText = \"83-22=\" },\n    @{ Row = 19; Col = 5; OldText = \"83-65=\"; NewText = \"62+23=\" },\n    @{ Row = 20; Col = 1; OldText = \"86-45=\"; NewText = \"32+2=\" },\n    @{ Row = 20; Col = 2; OldText = \"11+6=\"; NewText = \"69+17=\" },\n    @{ Row = 20; Col = 3; OldText = \"14+48=\"; NewText = \"49+44=\" },\n    @{ Row = 20; Col = 4; OldText = \"14+14=\"; NewText = \"76-56=\" },\n    @{ Row = 20; Col = 5; OldText = \"14+54=\"; NewText = \"58+25=\" }\n)\n\nforeach ($entry in $replacements) {\n    $cell = $t.Cell($entry.Row, $entry.Col)\n    # Cell.Range.Text includes the trailing end-of-cell marker (CR + cell\n    # mark), so strip it before comparing against the expected \"before\" text.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $entry.OldText) {\n        Write-Output (\"Warning: cell (\" + $entry.Row + \",\" + $entry.Col + \") was '\" + $current + \"', expected '\" + $entry.OldText + \"'\")\n    }\n    $cell.Range.Text = $entry.NewText\n}\n"}
